$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1091.1786
$ws.Range("J19").Value = 356.52942
$ws.Range("L19").Value = 356.52942
$ws.Range("N19").Value = -706.5294200000001

$ws.Range("H86").Value = 2264.8
$ws.Range("I86").Value = 1147.2142
$ws.Range("K86").Value = 1147.2142
$ws.Range("M86").Value = -24.21419999999989

$ws.Range("H89").Value = 2264.8
$ws.Range("I89").Value = 1147.2142
$ws.Range("K89").Value = 5736.071
$ws.Range("M89").Value = -120.0709999999999

$ws.Range("H92").Value = 872.9375
$ws.Range("I92").Value = 823.7857
$ws.Range("J92").Value = 1217
$ws.Range("K92").Value = 823.7857
$ws.Range("L92").Value = 1217
$ws.Range("M92").Value = 424.2143
$ws.Range("N92").Value = -3713

$ws.Range("H99").Value = 773.5
$ws.Range("I99").Value = 564.6667
$ws.Range("J99").Value = 1400
$ws.Range("K99").Value = 1694.0001
$ws.Range("L99").Value = 4200
$ws.Range("M99").Value = -196.0001
$ws.Range("N99").Value = -7196

$ws.Range("H100").Value = 1692.5
$ws.Range("J100").Value = 1550
$ws.Range("L100").Value = 1550
$ws.Range("N100").Value = -2632

$ws.Range("H132").Value = 2951.8057
$ws.Range("I132").Value = 3990.28
$ws.Range("J132").Value = 591.63635
$ws.Range("K132").Value = 11970.84
$ws.Range("L132").Value = 1774.90905
$ws.Range("M132").Value = -9440.84
$ws.Range("N132").Value = -6834.90905

$ws.Range("H135").Value = 13951179
$ws.Range("J135").Value = 35336476
$ws.Range("L135").Value = 318028284
$ws.Range("N135").Value = -318033354

$ws.Range("H137").Value = 35717704
$ws.Range("I137").Value = 1905.8334
$ws.Range("K137").Value = 5717.5002
$ws.Range("M137").Value = -3167.5002

$ws.Range("H138").Value = 4764.712
$ws.Range("I138").Value = 3330.2
$ws.Range("J138").Value = 5819.5
$ws.Range("K138").Value = 9990.599999999999
$ws.Range("L138").Value = 17458.5
$ws.Range("M138").Value = -4850.599999999999
$ws.Range("N138").Value = -27738.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26384
$ws.Range("I32").Value = 20043.732
$ws.Range("K32").Value = 20043.732
$ws.Range("M32").Value = -19756.732

$ws.Range("H45").Value = 1149.2941
$ws.Range("I45").Value = 1017.6
$ws.Range("J45").Value = 1204.1666
$ws.Range("K45").Value = 1017.6
$ws.Range("L45").Value = 1204.1666
$ws.Range("M45").Value = -640.6
$ws.Range("N45").Value = -1958.1666

$ws.Range("H74").Value = 10642648
$ws.Range("I74").Value = 18519386
$ws.Range("K74").Value = 18519386
$ws.Range("M74").Value = -18518512

$ws.Range("H77").Value = 10642648
$ws.Range("I77").Value = 18519386
$ws.Range("K77").Value = 92596930
$ws.Range("M77").Value = -92592562

$ws.Range("H110").Value = 1170
$ws.Range("I110").Value = 1070
$ws.Range("J110").Value = 1370
$ws.Range("K110").Value = 1070
$ws.Range("L110").Value = 1370
$ws.Range("M110").Value = 975
$ws.Range("N110").Value = -5460

$ws.Range("H139").Value = 71705
$ws.Range("J139").Value = 71705
$ws.Range("L139").Value = 71705
$ws.Range("N139").Value = -81985

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1470.3043
$ws.Range("I86").Value = 1337.7894
$ws.Range("K86").Value = 1337.7894
$ws.Range("M86").Value = -214.7893999999999

$ws.Range("H89").Value = 1470.3043
$ws.Range("I89").Value = 1337.7894
$ws.Range("K89").Value = 6688.946999999999
$ws.Range("M89").Value = -1072.946999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1693.4445
$ws.Range("I16").Value = 1480.1428
$ws.Range("J16").Value = 2440
$ws.Range("K16").Value = 1480.1428
$ws.Range("L16").Value = 2440
$ws.Range("M16").Value = -1193.1428
$ws.Range("N16").Value = -3014

$ws.Range("H31").Value = 17237.715
$ws.Range("I31").Value = 30857.883
$ws.Range("J31").Value = 7976
$ws.Range("K31").Value = 30857.883
$ws.Range("L31").Value = 7976
$ws.Range("M31").Value = -30562.883
$ws.Range("N31").Value = -8566

$ws.Range("H34").Value = 17237.715
$ws.Range("I34").Value = 30857.883
$ws.Range("J34").Value = 7976
$ws.Range("K34").Value = 30857.883
$ws.Range("L34").Value = 7976
$ws.Range("M34").Value = -30655.883
$ws.Range("N34").Value = -8380

$ws.Range("H93").Value = 18630.572
$ws.Range("I93").Value = 16735.666
$ws.Range("J93").Value = 30000
$ws.Range("K93").Value = 16735.666
$ws.Range("L93").Value = 30000
$ws.Range("M93").Value = -14863.666
$ws.Range("N93").Value = -33744

$ws.Range("H113").Value = 1693.4445
$ws.Range("I113").Value = 1480.1428
$ws.Range("J113").Value = 2440
$ws.Range("K113").Value = 1480.1428
$ws.Range("L113").Value = 2440
$ws.Range("M113").Value = 689.8571999999999
$ws.Range("N113").Value = -6780

$ws.Range("H141").Value = 55803.145
$ws.Range("J141").Value = 70865.2
$ws.Range("L141").Value = 70865.2
$ws.Range("N141").Value = -81225.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 15625725
$ws.Range("I5").Value = 687.7083
$ws.Range("J5").Value = 62500836
$ws.Range("K5").Value = 2063.1249
$ws.Range("L5").Value = 187502508
$ws.Range("M5").Value = -1951.1249
$ws.Range("N5").Value = -187502732

$ws.Range("H10").Value = 158.375
$ws.Range("I10").Value = 137.85715
$ws.Range("K10").Value = 413.57145
$ws.Range("M10").Value = -274.57145

$ws.Range("H107").Value = 80250.03999999999
$ws.Range("I107").Value = 59069.293
$ws.Range("J107").Value = 125259.125
$ws.Range("K107").Value = 177207.879
$ws.Range("L107").Value = 375777.375
$ws.Range("M107").Value = -175287.879
$ws.Range("N107").Value = -379617.375

$ws.Range("H113").Value = 547.8
$ws.Range("I113").Value = 516.4
$ws.Range("J113").Value = 579.2
$ws.Range("K113").Value = 1549.2
$ws.Range("L113").Value = 1737.6
$ws.Range("M113").Value = 620.8000000000002
$ws.Range("N113").Value = -6077.6

$ws.Range("H125").Value = 2931.6667
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 3318
$ws.Range("K125").Value = 3000
$ws.Range("L125").Value = 9954
$ws.Range("M125").Value = 1920
$ws.Range("N125").Value = -19794

$ws.Range("H131").Value = 37638560
$ws.Range("I131").Value = 74080970
$ws.Range("J131").Value = 22730302
$ws.Range("K131").Value = 222242910
$ws.Range("L131").Value = 68190906
$ws.Range("M131").Value = -222237870
$ws.Range("N131").Value = -68200986

$ws.Range("H132").Value = 31250978
$ws.Range("I132").Value = 41667640
$ws.Range("J132").Value = 996.875
$ws.Range("K132").Value = 375008760
$ws.Range("L132").Value = 8971.875
$ws.Range("M132").Value = -375006230
$ws.Range("N132").Value = -14031.875

$ws.Range("H135").Value = 15625725
$ws.Range("I135").Value = 687.7083
$ws.Range("J135").Value = 62500836
$ws.Range("K135").Value = 6189.3747
$ws.Range("L135").Value = 562507524
$ws.Range("M135").Value = -3654.3747
$ws.Range("N135").Value = -562512594

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5954.923
$ws.Range("I122").Value = 7639.25
$ws.Range("J122").Value = 3260
$ws.Range("K122").Value = 22917.75
$ws.Range("L122").Value = 9780
$ws.Range("M122").Value = -20467.75
$ws.Range("N122").Value = -14680

$ws.Range("H132").Value = 2253.5715
$ws.Range("I132").Value = 1728
$ws.Range("J132").Value = 3199.6
$ws.Range("K132").Value = 5184
$ws.Range("L132").Value = 9598.799999999999
$ws.Range("M132").Value = -2654
$ws.Range("N132").Value = -14658.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8786.817999999999
$ws.Range("I132").Value = 9522.611000000001
$ws.Range("J132").Value = 5475.75
$ws.Range("K132").Value = 28567.833
$ws.Range("L132").Value = 16427.25
$ws.Range("M132").Value = -26037.833
$ws.Range("N132").Value = -21487.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2384.913
$ws.Range("I126").Value = 2113.3845
$ws.Range("J126").Value = 2737.9
$ws.Range("K126").Value = 6340.1535
$ws.Range("L126").Value = 8213.700000000001
$ws.Range("M126").Value = -3870.1535
$ws.Range("N126").Value = -13153.7

$ws.Range("H136").Value = 22556.834
$ws.Range("I136").Value = 26003.2
$ws.Range("J136").Value = 5325
$ws.Range("K136").Value = 78009.60000000001
$ws.Range("L136").Value = 15975
$ws.Range("M136").Value = -75459.60000000001
$ws.Range("N136").Value = -21075
